$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.008803129196167
$ws.Range("B1").Value = 1.415010213851929
$ws.Range("C1").Value = 3.750451803207397
$ws.Range("D1").Value = 2.264014720916748
$ws.Range("E1").Value = 0.7641845345497131
